$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45412.76131387657

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 15).Value = $newDate
}
